$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Shift the tail of the sheet (old rows 34-45) down by 4 rows (-> 38-49)
#    without disturbing their existing formatting. We do this by inserting a
#    copy of a simple, uniformly-formatted 4-row block (42:45, which only
#    ever carries style s=14 on column B) at row 34; Excel's "insert copied
#    cells" shifts everything at/after row 34 down by the size of the copied
#    block while keeping the original rows' own formatting intact.
# ---------------------------------------------------------------------------
$ws.Rows("42:45").Copy()
$ws.Rows("34:37").Insert(-4121)
$ws.Application.CutCopyMode = $false

# The freshly inserted rows 34-37 picked up incidental formatting from the
# insert heuristic; wipe them so we can build the new content cleanly.
$ws.Range("A34:B37").ClearContents()
$ws.Range("A34:B37").ClearFormats()

# ---------------------------------------------------------------------------
# 2. Row 34: new "LANGKAH-LANGKAH PEMBELAJARAN" section marker row. This
#    mirrors the existing section-marker rows (e.g. row 16 "SARANA DAN
#    PRASARANA"): column A bold-on-yellow with the section title, column B
#    bold-on-yellow with the literal text "DESKRIPSI".
# ---------------------------------------------------------------------------
$ws.Range("A16:B16").Copy()
$ws.Range("A34:B34").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("A34").Value = "LANGKAH-LANGKAH PEMBELAJARAN"
$ws.Range("B34").Value = "DESKRIPSI"

# ---------------------------------------------------------------------------
# 3. Rows 35-37: "Kegiatan Awal" / "Kegiatan Inti" / "Kegiatan Penutup".
#    New bold/Calibri/12pt font on a bordered, vertically centered,
#    word-wrapped cell (closest existing style is row 22's label cell, which
#    has the border + vertical-center but isn't bold and doesn't wrap).
# ---------------------------------------------------------------------------
$ws.Range("A22:B22").Copy()
$ws.Range("A35:B35").PasteSpecial(-4122)
$ws.Range("A36:B36").PasteSpecial(-4122)
$ws.Range("A37:B37").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("A35").Value = "Kegiatan Awal"
$ws.Range("A36").Value = "Kegiatan Inti"
$ws.Range("A37").Value = "Kegiatan Penutup"

$labelRange = $ws.Range("A35:A37")
$labelRange.Font.Bold = $true
$labelRange.Font.Name = "Calibri"
$labelRange.Font.Size = 12
$labelRange.Font.Color = 0
$labelRange.WrapText = $true
$labelRange.VerticalAlignment = -4108

$ws.Range("B35:B37").Value = ""

$ws.Rows(35).RowHeight = 56.25
$ws.Rows(36).RowHeight = 52.5
$ws.Rows(37).RowHeight = 49.5

# Row 38 stays a plain spacer row (same as the other B-only s=14 rows).
$ws.Range("B38").Value = ""

# ---------------------------------------------------------------------------
# 4. Column A width grew to fit the new labels; bring the sheet dimension
#    and view in line with the new layout.
# ---------------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 39.8

$ws.Range("A34").Select()
$ws.Application.ActiveWindow.ScrollRow = 34
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("D37").Select()

$ws.PageSetup.Orientation = 1
